$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New row 9: Vorbereiten Präsentationstechnikkurs (copy formatting from row 11, which
# already has the right date / text / number styles)
$ws.Range("A11:C11").Copy($ws.Range("A9:C9"))
$ws.Cells.Item(9, 1).Value = 41399
$ws.Cells.Item(9, 2).Value = "Vorbereiten Präsentationstechnikkurs"
$ws.Cells.Item(9, 3).Value = 1

# New row 10: Präsentationstechnikkurs (duplicate of the old row 11 content)
$ws.Range("A11:C11").Copy($ws.Range("A10:C10"))
$ws.Cells.Item(10, 1).Value = 41400
$ws.Cells.Item(10, 2).Value = "Präsentationstechnikkurs"
$ws.Cells.Item(10, 3).Value = 8

# Row 11 now becomes: Feedback aus Kickoff verarbeiten
$ws.Cells.Item(11, 1).Value = 41400
$ws.Cells.Item(11, 2).Value = "Feedback aus Kickoff verarbeiten"
$ws.Cells.Item(11, 3).Value = 3

# Move the TOTAL row up from row 35 to row 34 (copy formatting, then fix values/formula)
$ws.Range("A35:C35").Copy($ws.Range("A34:C34"))
$ws.Cells.Item(34, 1).ClearContents()
$ws.Cells.Item(34, 2).Value = "TOTAL"
$ws.Cells.Item(34, 3).Formula = "=SUM(C2:C33)"

# Remove the old total row
$ws.Rows.Item(35).Delete()

$ws.Range("C12").Select()

$wb.Save()
